$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-14 Thursday" "2024-11-15 Friday"

Replace-Text "793÷8=99, 1" "796÷7=113, 5"
Replace-Text "878÷2=439, 0" "713÷3=237, 2"
Replace-Text "889÷5=177, 4" "512÷3=170, 2"
Replace-Text "518÷9=57, 5" "195÷5=39, 0"
Replace-Text "432÷5=86, 2" "899÷6=149, 5"

Replace-Text "648÷8=81, 0" "898÷7=128, 2"
Replace-Text "549÷8=68, 5" "199÷5=39, 4"
Replace-Text "114÷4=28, 2" "711÷7=101, 4"
Replace-Text "120÷6=20, 0" "499÷2=249, 1"
Replace-Text "843÷7=120, 3" "950÷8=118, 6"

Replace-Text "426÷6=71, 0" "674÷9=74, 8"
Replace-Text "726÷6=121, 0" "971÷6=161, 5"
Replace-Text "637÷2=318, 1" "445÷7=63, 4"
Replace-Text "492÷4=123, 0" "225÷2=112, 1"
Replace-Text "105÷4=26, 1" "190÷8=23, 6"

Replace-Text "307÷6=51, 1" "277÷9=30, 7"
Replace-Text "612÷3=204, 0" "883÷3=294, 1"
Replace-Text "457÷9=50, 7" "517÷4=129, 1"
Replace-Text "464÷6=77, 2" "800÷6=133, 2"
Replace-Text "213÷8=26, 5" "516÷6=86, 0"

Replace-Text "643÷5=128, 3" "365÷2=182, 1"
Replace-Text "632÷6=105, 2" "922÷9=102, 4"
Replace-Text "728÷9=80, 8" "601÷5=120, 1"
Replace-Text "861÷6=143, 3" "606÷5=121, 1"
Replace-Text "304÷6=50, 4" "201÷8=25, 1"
